# [MOSIP-14369] Fix: boolean values
#
# Column I ("is_active") currently stores the boolean result of the
# formula =TRUE() (numeric cell, value 1). The fix replaces those cells
# with the literal text "TRUE" instead of a computed boolean, for rows
# 2 through 11.
#
# NOTE: Directly assigning the string "TRUE" via .Value/.Value2/.Formula
# is auto-recognised by Excel as the Boolean TRUE, which would just
# reproduce the original (undesired) boolean cell. To force a literal
# text cell, we first compute the text via TEXT(TRUE(),"General") and
# then flatten the formula down to a plain value with Copy/PasteSpecial
# (values only), which keeps the existing cell style/number format
# intact and yields a real text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Formula = '=TEXT(TRUE(),"General")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = $false

# Match the updated selection left behind in the sheet (I2:I11 selected,
# active cell I2) instead of the original I:I / I1 selection.
$null = $ws.Range("I2:I11").Select()
